# Updates cryptos list (Coin / Link / Price / Volume(1h)) to the latest
# scraped values. Mirrors the source GitHub Actions refresh job.
#
# The "Price" column stores numeric-looking values (e.g. "583.92",
# "71.229.93") as TEXT so that thousand-dot formatting and trailing
# zeros survive untouched -- Excel would otherwise silently coerce a
# plain decimal string into a real Number. Set-TextValue forces the
# text interpretation the same way a user would (leading apostrophe)
# and then restores the "Normal" style so no stray NumberFormat/
# quote-prefix formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "71.229.93"
$ws.Range("E2").Value = "  +3.01%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.577.26"
$ws.Range("E3").Value = "  +1.31%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.14%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "583.92"
$ws.Range("E5").Value = "  +2.23%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "186.82"
$ws.Range("E6").Value = "  +2.69%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.567.44"
$ws.Range("E7").Value = "  +1.25%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +1.36%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  -0.07%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.217"
$ws.Range("E10").Value = "  +14.50%  "

# Row 11 - Cardano
Set-TextValue $ws.Range("D11") "0.655"
$ws.Range("E11").Value = "  +2.78%  "

# Row 12 - Avalanche
Set-TextValue $ws.Range("D12") "54.67"
$ws.Range("E12").Value = "  +1.93%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +5.78%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "9.56"
$ws.Range("E14").Value = "  +1.08%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.057.15"
$ws.Range("E15").Value = "  -1.08%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "71.135.41"
$ws.Range("E16").Value = "  +2.77%  "

# Row 17 - Chainlink
Set-TextValue $ws.Range("D17") "19.28"
$ws.Range("E17").Value = "  +0.10%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.576.19"
$ws.Range("E18").Value = "  +1.18%  "

# Row 19 - Uniswap
Set-TextValue $ws.Range("D19") "12.40"
$ws.Range("E19").Value = "  -0.50%  "

# Row 20 - BitcoinCash
Set-TextValue $ws.Range("D20") "566.46"
$ws.Range("E20").Value = "  +5.31%  "

# Row 21 - TRON
$ws.Range("E21").Value = "  +0.63%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  -1.68%  "

# Row 23 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D23") "17.56"
$ws.Range("E23").Value = "  -11.26%  "

# Row 24 - Toncoin
Set-TextValue $ws.Range("D24") "5.07"
$ws.Range("E24").Value = "  +2.80%  "

# Row 25 - PancakeSwap
Set-TextValue $ws.Range("D25") "4.59"
$ws.Range("E25").Value = "  +4.92%  "

# Row 26 - Litecoin
Set-TextValue $ws.Range("D26") "94.75"
$ws.Range("E26").Value = "  +0.62%  "

# Row 27 - RenderToken
Set-TextValue $ws.Range("D27") "11.31"
$ws.Range("E27").Value = "  +2.48%  "

# Row 28 - ImmutableX
$ws.Range("E28").Value = "  +1.47%  "

# Row 29 - Filecoin
$ws.Range("E29").Value = "  +1.35%  "

# Row 30 - EthereumClassic
Set-TextValue $ws.Range("D30") "32.68"
$ws.Range("E30").Value = "  +2.94%  "

# Row 31 - NEARProtocol
Set-TextValue $ws.Range("D31") "7.27"
$ws.Range("E31").Value = "  -0.69%  "

# Row 32 - Cosmos
Set-TextValue $ws.Range("D32") "12.37"
$ws.Range("E32").Value = "  -1.47%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  +1.16%  "

# Row 34 - OKB
Set-TextValue $ws.Range("D34") "64.25"
$ws.Range("E34").Value = "  -1.04%  "

# Row 35 - Fetch.AI
Set-TextValue $ws.Range("D35") "3.39"
$ws.Range("E35").Value = "  +8.09%  "

# Row 36 - Bittensor
Set-TextValue $ws.Range("D36") "553.78"
$ws.Range("E36").Value = "  -3.23%  "

# Row 37 - TheGraph
Set-TextValue $ws.Range("D37") "0.420"
$ws.Range("E37").Value = "  +5.47%  "

# Row 38 - PEPE
$ws.Range("D38").Value = "0.0₃0809"
$ws.Range("E38").Value = "  +6.13%  "

# Row 39 - InjectiveProtocol
Set-TextValue $ws.Range("D39") "37.73"
$ws.Range("E39").Value = "  -1.24%  "

# Row 40 - Dai
$ws.Range("E40").Value = "  +0.10%  "

# Row 41 - dogwifhat
$ws.Range("E41").Value = "  +6.03%  "

# Row 42 - Maker
$ws.Range("D42").Value = "3.504.51"
$ws.Range("E42").Value = "  +10.63%  "

# Row 43 - Stacks
$ws.Range("E43").Value = "  +2.46%  "

# Row 44 - Kaspa
$ws.Range("E44").Value = "  +2.44%  "

# Row 45 - VeChain
Set-TextValue $ws.Range("D45") "0.0448"
$ws.Range("E45").Value = "  +1.18%  "

# Row 46 - ThetaToken
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D46") "2.95"
$ws.Range("E46").Value = "  -0.32%  "

# Row 47 - ApeXProtocol
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D47") "3.46"
$ws.Range("E47").Value = "  -1.16%  "

# Row 48
Set-TextValue $ws.Range("D48") "9.37"
$ws.Range("E48").Value = "  +1.52%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.138"
$ws.Range("E49").Value = "  +2.92%  "

# Row 50 - OceanProtocol
$ws.Range("B50").Value = "OceanProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
Set-TextValue $ws.Range("D50") "1.48"
$ws.Range("E50").Value = "  +8.29%  "

# Row 51 - FirstDigitalUSD
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D51") "0.998"
$ws.Range("E51").Value = "  -0.32%  "
